$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 8 ("Incorrect Username and Incorrect Password" test
# case, whose B8 hyperlink had no display text) -- this shifts rows
# 9-11 up to become rows 8-10.
$ws.Rows("8").Delete()

# The hyperlink ranges don't auto-shift with the row delete, so rebuild
# the hyperlinks collection to match the new layout: the surviving
# "prasoona@testcase.com" mail links now sit on B8 (was B9) and B6
# (unchanged).
$ws.Hyperlinks.Delete()

$hl1 = $ws.Hyperlinks.Add($ws.Range("B8"), "mailto:prasoona@testcase.com")
$hl1.TextToDisplay = "prasoona@testcase.com"
$ws.Range("B8").Value2 = "alphatestcase2020@gmail.com"
$ws.Range("B8").Font.FontStyle = $ws.Range("B8").Font.FontStyle

$hl2 = $ws.Hyperlinks.Add($ws.Range("B6"), "mailto:prasoona@testcase.com")
$hl2.TextToDisplay = "prasoona@testcase.com"
$ws.Range("B6").Value2 = "alphatestcase2020@gmail.com"
$ws.Range("B6").Font.FontStyle = $ws.Range("B6").Font.FontStyle

# Move the selection to reflect where the user ended up after the edit.
$ws.Range("A10").Select()
